# Auto-generated edit script: apply scheduled market-price refresh to Hyperion_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 841.5833
$ws.Range("I4").Value = 212.5
$ws.Range("K4").Value = 212.5
$ws.Range("M4").Value = -98.5
$ws.Range("H33").Value = 569.2308
$ws.Range("I33").Value = 610.7826
$ws.Range("J33").Value = 250.66667
$ws.Range("K33").Value = 610.7826
$ws.Range("L33").Value = 250.66667
$ws.Range("M33").Value = -381.7826
$ws.Range("N33").Value = -708.6666700000001
$ws.Range("H64").Value = 5271.2856
$ws.Range("I64").Value = 3200
$ws.Range("J64").Value = 6099.8
$ws.Range("K64").Value = 3200
$ws.Range("L64").Value = 6099.8
$ws.Range("M64").Value = -2952
$ws.Range("N64").Value = -6595.8
$ws.Range("H67").Value = 5271.2856
$ws.Range("I67").Value = 3200
$ws.Range("J67").Value = 6099.8
$ws.Range("K67").Value = 3200
$ws.Range("L67").Value = 6099.8
$ws.Range("M67").Value = -2342
$ws.Range("N67").Value = -7815.8
$ws.Range("H88").Value = 6859.6665
$ws.Range("I88").Value = 5277.5
$ws.Range("J88").Value = 7311.7144
$ws.Range("K88").Value = 5277.5
$ws.Range("L88").Value = 7311.7144
$ws.Range("M88").Value = -4871.5
$ws.Range("N88").Value = -8123.7144
$ws.Range("H91").Value = 6859.6665
$ws.Range("I91").Value = 5277.5
$ws.Range("J91").Value = 7311.7144
$ws.Range("K91").Value = 5277.5
$ws.Range("L91").Value = 7311.7144
$ws.Range("M91").Value = -3873.5
$ws.Range("N91").Value = -10119.7144
$ws.Range("H96").Value = 153093.47
$ws.Range("I96").Value = 541.46155
$ws.Range("J96").Value = 483622.84
$ws.Range("K96").Value = 1624.38465
$ws.Range("L96").Value = 1450868.52
$ws.Range("M96").Value = -251.38465
$ws.Range("N96").Value = -1453614.52
$ws.Range("H111").Value = 374.15
$ws.Range("I111").Value = 372.78946
$ws.Range("J111").Value = 400
$ws.Range("K111").Value = 1118.36838
$ws.Range("L111").Value = 1200
$ws.Range("M111").Value = 1948.63162
$ws.Range("N111").Value = -7334
$ws.Range("H113").Value = 8528.714
$ws.Range("I113").Value = 11158.75
$ws.Range("J113").Value = 7476.7
$ws.Range("K113").Value = 11158.75
$ws.Range("L113").Value = 7476.7
$ws.Range("M113").Value = -7904.75
$ws.Range("N113").Value = -13984.7
$ws.Range("H115").Value = 399.8889
$ws.Range("I115").Value = 402
$ws.Range("K115").Value = 1206
$ws.Range("M115").Value = 361
$ws.Range("H116").Value = 5534.05
$ws.Range("I116").Value = 4875
$ws.Range("J116").Value = 6073.273
$ws.Range("K116").Value = 4875
$ws.Range("L116").Value = 6073.273
$ws.Range("M116").Value = -1433
$ws.Range("N116").Value = -12957.273
$ws.Range("H132").Value = 4372.638
$ws.Range("I132").Value = 4432.1816
$ws.Range("K132").Value = 13296.5448
$ws.Range("M132").Value = -10766.5448
$ws.Range("H135").Value = 581.7273
$ws.Range("I135").Value = 489.9
$ws.Range("J135").Value = 1500
$ws.Range("K135").Value = 4409.099999999999
$ws.Range("L135").Value = 13500
$ws.Range("M135").Value = -1874.099999999999
$ws.Range("N135").Value = -18570
$ws.Range("H137").Value = 178811.3
$ws.Range("I137").Value = 198537.11
$ws.Range("J137").Value = 1279
$ws.Range("K137").Value = 595611.33
$ws.Range("L137").Value = 3837
$ws.Range("M137").Value = -593061.33
$ws.Range("N137").Value = -8937
$ws.Range("H138").Value = 4084.6145
$ws.Range("I138").Value = 3759
$ws.Range("J138").Value = 4156.4414
$ws.Range("K138").Value = 11277
$ws.Range("L138").Value = 12469.3242
$ws.Range("M138").Value = -6137
$ws.Range("N138").Value = -22749.3242

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 106018.6
$ws.Range("I45").Value = 170632.17
$ws.Range("K45").Value = 170632.17
$ws.Range("M45").Value = -170255.17
$ws.Range("H74").Value = 53983.12
$ws.Range("I74").Value = 11674.556
$ws.Range("K74").Value = 11674.556
$ws.Range("M74").Value = -10800.556
$ws.Range("H77").Value = 53983.12
$ws.Range("I77").Value = 11674.556
$ws.Range("K77").Value = 58372.78
$ws.Range("M77").Value = -54004.78
$ws.Range("H122").Value = 4009.5833
$ws.Range("I122").Value = 3442.5334
$ws.Range("J122").Value = 4954.6665
$ws.Range("K122").Value = 10327.6002
$ws.Range("L122").Value = 14863.9995
$ws.Range("M122").Value = -7877.600199999999
$ws.Range("N122").Value = -19763.9995
$ws.Range("H132").Value = 18000.25
$ws.Range("J132").Value = 7599
$ws.Range("L132").Value = 22797
$ws.Range("N132").Value = -27857
$ws.Range("H140").Value = 79321.25
$ws.Range("J140").Value = 79321.25
$ws.Range("L140").Value = 79321.25
$ws.Range("N140").Value = -89681.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3989.9092
$ws.Range("I20").Value = 3565.3125
$ws.Range("K20").Value = 3565.3125
$ws.Range("M20").Value = -3318.3125
$ws.Range("H33").Value = 31333.334
$ws.Range("I33").Value = 4000
$ws.Range("K33").Value = 4000
$ws.Range("M33").Value = -3664
$ws.Range("H86").Value = 5336.2144
$ws.Range("I86").Value = 5795.7427
$ws.Range("J86").Value = 3038.5715
$ws.Range("K86").Value = 5795.7427
$ws.Range("L86").Value = 3038.5715
$ws.Range("M86").Value = -4672.7427
$ws.Range("N86").Value = -5284.5715
$ws.Range("H89").Value = 5336.2144
$ws.Range("I89").Value = 5795.7427
$ws.Range("J89").Value = 3038.5715
$ws.Range("K89").Value = 28978.7135
$ws.Range("L89").Value = 15192.8575
$ws.Range("M89").Value = -23362.7135
$ws.Range("N89").Value = -26424.8575
$ws.Range("H105").Value = 113437.336
$ws.Range("I105").Value = 144992.28
$ws.Range("J105").Value = 2995
$ws.Range("K105").Value = 144992.28
$ws.Range("L105").Value = 2995
$ws.Range("M105").Value = -143245.28
$ws.Range("N105").Value = -6489

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2458.5386
$ws.Range("I16").Value = 2343
$ws.Range("K16").Value = 2343
$ws.Range("M16").Value = -2056
$ws.Range("H31").Value = 34999.516
$ws.Range("I31").Value = 8285.733
$ws.Range("K31").Value = 8285.733
$ws.Range("M31").Value = -7990.733
$ws.Range("H34").Value = 34999.516
$ws.Range("I34").Value = 8285.733
$ws.Range("K34").Value = 8285.733
$ws.Range("M34").Value = -8083.733
$ws.Range("H113").Value = 2458.5386
$ws.Range("I113").Value = 2343
$ws.Range("K113").Value = 2343
$ws.Range("M113").Value = -173
$ws.Range("H122").Value = 4102.6
$ws.Range("I122").Value = 3212.4
$ws.Range("K122").Value = 9637.200000000001
$ws.Range("M122").Value = -7187.200000000001
$ws.Range("H132").Value = 57504.95
$ws.Range("J132").Value = 10000
$ws.Range("L132").Value = 30000
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1271.3334
$ws.Range("I60").Value = 1146.6471
$ws.Range("K60").Value = 3439.9413
$ws.Range("M60").Value = -3188.9413
$ws.Range("H68").Value = 3999.5454
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 3999.5454
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 11998.6362
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -13620.6362
$ws.Range("H71").Value = 3999.5454
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 3999.5454
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 35995.9086
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -44107.9086
$ws.Range("H109").Value = 3657.8
$ws.Range("J109").Value = 3068.5715
$ws.Range("L109").Value = 9205.7145
$ws.Range("N109").Value = -11285.7145
$ws.Range("H136").Value = 2010.5264
$ws.Range("I136").Value = 1381.8182
$ws.Range("K136").Value = 4145.4546
$ws.Range("M136").Value = 954.5454
$ws.Range("H138").Value = 24776.846
$ws.Range("J138").Value = 31860
$ws.Range("L138").Value = 95580
$ws.Range("N138").Value = -105860
$ws.Range("H140").Value = 2261.257
$ws.Range("I140").Value = 1313.2307
$ws.Range("K140").Value = 3939.6921
$ws.Range("M140").Value = 1240.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 26723
$ws.Range("J96").Value = 26723
$ws.Range("L96").Value = 26723
$ws.Range("N96").Value = -32215
$ws.Range("H102").Value = 39915.15
$ws.Range("I102").Value = 2083.5293
$ws.Range("K102").Value = 2083.5293
$ws.Range("M102").Value = -461.5293000000001
$ws.Range("H106").Value = 33000
$ws.Range("J106").Value = 33000
$ws.Range("L106").Value = 33000
$ws.Range("N106").Value = -35524
$ws.Range("H122").Value = 991998.25
$ws.Range("I122").Value = 1273497.8
$ws.Range("J122").Value = 6750
$ws.Range("K122").Value = 3820493.4
$ws.Range("L122").Value = 20250
$ws.Range("M122").Value = -3818043.4
$ws.Range("N122").Value = -25150
$ws.Range("H126").Value = 4262.8
$ws.Range("I126").Value = 3424
$ws.Range("J126").Value = 4996.75
$ws.Range("K126").Value = 10272
$ws.Range("L126").Value = 14990.25
$ws.Range("M126").Value = -7802
$ws.Range("N126").Value = -19930.25
$ws.Range("H132").Value = 8655.652
$ws.Range("I132").Value = 3651.9285
$ws.Range("J132").Value = 16439.223
$ws.Range("K132").Value = 10955.7855
$ws.Range("L132").Value = 49317.66900000001
$ws.Range("M132").Value = -8425.7855
$ws.Range("N132").Value = -54377.66900000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5008151
$ws.Range("I40").Value = 7148087.5
$ws.Range("K40").Value = 7148087.5
$ws.Range("M40").Value = -7147951.5
$ws.Range("H46").Value = 6487.8823
$ws.Range("I46").Value = 3885.5715
$ws.Range("J46").Value = 8309.5
$ws.Range("K46").Value = 3885.5715
$ws.Range("L46").Value = 8309.5
$ws.Range("M46").Value = -3697.5715
$ws.Range("N46").Value = -8685.5
$ws.Range("H68").Value = 5087.5
$ws.Range("I68").Value = 3177.6
$ws.Range("J68").Value = 6997.4
$ws.Range("K68").Value = 3177.6
$ws.Range("L68").Value = 6997.4
$ws.Range("M68").Value = -2428.6
$ws.Range("N68").Value = -8495.4
$ws.Range("H71").Value = 5087.5
$ws.Range("I71").Value = 3177.6
$ws.Range("J71").Value = 6997.4
$ws.Range("K71").Value = 15888
$ws.Range("L71").Value = 34987
$ws.Range("M71").Value = -12144
$ws.Range("N71").Value = -42475
$ws.Range("H100").Value = 4888.778
$ws.Range("I100").Value = 3800
$ws.Range("K100").Value = 3800
$ws.Range("M100").Value = -3259
$ws.Range("H122").Value = 8110.4287
$ws.Range("I122").Value = 5355
$ws.Range("K122").Value = 16065
$ws.Range("M122").Value = -13615
$ws.Range("H132").Value = 2245.4
$ws.Range("I132").Value = 2064.8572
$ws.Range("J132").Value = 2666.6667
$ws.Range("K132").Value = 6194.571599999999
$ws.Range("L132").Value = 8000.000100000001
$ws.Range("M132").Value = -3664.571599999999
$ws.Range("N132").Value = -13060.0001
$ws.Range("H136").Value = 22258.98
$ws.Range("I136").Value = 35070.266
$ws.Range("J136").Value = 3957.1428
$ws.Range("K136").Value = 105210.798
$ws.Range("L136").Value = 11871.4284
$ws.Range("M136").Value = -102660.798
$ws.Range("N136").Value = -16971.4284
$ws.Range("H138").Value = 125000
$ws.Range("J138").Value = 125000
$ws.Range("L138").Value = 125000
$ws.Range("N138").Value = -135280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 12230.143
$ws.Range("J32").Value = 16521.75
$ws.Range("L32").Value = 16521.75
$ws.Range("N32").Value = -17155.75
$ws.Range("H39").Value = 1709
$ws.Range("I39").Value = 1709
$ws.Range("K39").Value = 1709
$ws.Range("M39").Value = -1296
$ws.Range("H81").Value = 625
$ws.Range("I81").Value = 666.6667
$ws.Range("K81").Value = 1333.3334
$ws.Range("M81").Value = -272.3334
$ws.Range("H84").Value = 625
$ws.Range("I84").Value = 666.6667
$ws.Range("K84").Value = 6666.666999999999
$ws.Range("M84").Value = -1362.666999999999
$ws.Range("H105").Value = 78663
$ws.Range("J105").Value = 78663
$ws.Range("L105").Value = 78663
$ws.Range("N105").Value = -85651
$ws.Range("H107").Value = 2691
$ws.Range("I107").Value = 2961.4
$ws.Range("J107").Value = 663
$ws.Range("K107").Value = 8884.200000000001
$ws.Range("L107").Value = 1989
$ws.Range("M107").Value = -6964.200000000001
$ws.Range("N107").Value = -5829
$ws.Range("H122").Value = 5287.9
$ws.Range("I122").Value = 4814
$ws.Range("K122").Value = 14442
$ws.Range("M122").Value = -11992
$ws.Range("H126").Value = 3535
$ws.Range("I126").Value = 3227.647
$ws.Range("K126").Value = 9682.940999999999
$ws.Range("M126").Value = -7212.940999999999
$ws.Range("H132").Value = 373563.72
$ws.Range("I132").Value = 8676.65
$ws.Range("K132").Value = 26029.95
$ws.Range("M132").Value = -23499.95
$ws.Range("H136").Value = 5598.6924
$ws.Range("I136").Value = 6400.067
$ws.Range("K136").Value = 19200.201
$ws.Range("M136").Value = -16650.201
